$wb = $excel.ActiveWorkbook

# --- Rename sheet2 DATA -> RETAIL_DATA ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "RETAIL_DATA"

# =========================================================================
# Sheet1 (RUNMANAGER)
# =========================================================================
$ws1.Range("B2").Value = "To check whether the user can successfully login and logout"
$ws1.Range("C2").Value = "yes"
$ws1.Range("D2").Value = "'1"
$ws1.Range("E2").Value = "'1"

$ws1.Range("A3").Value = "newTest"
$ws1.Range("B3").Value = "To check this test runs"
$ws1.Range("C3").Value = "yes"
$ws1.Range("D3").Value = "'1"
$ws1.Range("E3").Value = "'1"

# =========================================================================
# Sheet2 (RETAIL_DATA)
# =========================================================================
$ws2.Range("C1").Value = "browser"
$ws2.Range("D1").Value = "username"
$ws2.Range("E1").Value = "password"

$ws2.Range("B2").Value = "yes"
$ws2.Range("C2").Value = "chrome"
$ws2.Range("D2").Value = "spcbtest"
$ws2.Range("E2").Value = "Asdf@123"

$ws2.Range("A3").Value = "newTest"
$ws2.Range("B3").Value = "yes"
$ws2.Range("C3").Value = "chrome"
$ws2.Range("D3").Value = "spcb"
$ws2.Range("E3").Value = "Asdf@123"

# Remove the old data from rows 4-6 (A:D), keep column E cells (styled) for
# the hyperlink-only rows, then rebuild hyperlinks to just E2/E3.
$ws2.Range("A4:D6").ClearContents()
$ws2.Range("E4:E6").ClearContents()

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("E2"), "mailto:Asdf@123")
$ws2.Range("E2").Style = "Hyperlink"
$ws2.Hyperlinks.Add($ws2.Range("E3"), "mailto:Asdf@123")
$ws2.Range("E3").Style = "Hyperlink"

# New row 7: only E7, carrying the same (Hyperlink) style as E6, no value/link.
$ws2.Range("E6").Copy()
$ws2.Range("E7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# =========================================================================
# Selections / active sheet
# =========================================================================
$ws2.Range("C20").Select()
$ws1.Range("C10").Select()
$ws1.Activate()
